$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix typo in author name
$ws.Range("A2").Value = "Kurt Tuohy and Michael Miller"

# 2. Fix "Power & ctrl" text (remove stray backslash before ampersand)
$ws.Range("A12").Value = "Power & ctrl   "

# 3. Add new footnote describing the table, in cell A18
$ws.Range("A18").Value = "For each model type, displays Spearman rank-order correlation between ground-truth and predicted correspondence of utterances to schemas."

# 4. Move the active selection to reflect where the user finished editing (below the new note)
$ws.Range("A19").Select()

# 5. Slightly widen column A so the longer note/author text is not clipped
$ws.Columns("A").ColumnWidth = 15.5
